$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo "Ngnix 3" -> "Nginx 3"
$ws.Range("A5").Value = "Nginx 3"

# Trim trailing space "Apache 3 " -> "Apache 3"
$ws.Range("A11").Value = "Apache 3"

# Widen column A
$ws.Columns.Item(1).ColumnWidth = 25.86
